$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Work in Progress")

$ws.Range("A4").Value = "The Long-run Effects of Unconditional Cash Transfers: Evidence from the Kenya General Equilibrium Study"
$ws.Range("B4").Value = "with Dennis Egger, Edward Miguel and Michael Walker"
$ws.Range("C4").Value = "Recent studies document positive short-run effects of unconditional cash transfers (UCTs) on cash recipients and spillover effects on non-recipients. But modest sample sizes and challenges with tracking households over time have limited research on the long-run effects of UCTs. We study the long-run effects of the Kenya General Equilibrium Study (KGES) on recipient households and the local economy using census and survey data collected up to ten years post-transfers. An “endline 2” completed 4-7 years after the transfers collected census data from each household and enterprise in the study area, and a representative survey of more than 10,000 households and firms obtained more detailed consumption and production measures for a subset, including from those that migrated out of the study area. We collected this same information in an “endline 3” completed 7-10 years after the experimental start. Preliminary results from endline 2 show persistent consumption gains among recipients and expansion of non-farm enterprise revenue. Ongoing analysis examines whether these effects persisted to endline 3 and aims to estimate long-run transfer multipliers."

$ws.Range("A4:C4").WrapText = $true

$ws.Activate()
$ws.Range("A3").Select()
$ws.Range("B5").Select()
